$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Drashti"
$ws.Range("B2").Value = "Raja"

$ws.Range("B2").Select()
